$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "66.030.93"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  +1.65%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.210.97"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +1.28%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "602.24"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +3.90%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "154.00"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +1.69%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.211.22"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +1.29%  "
$ws.Range("E9").Value = "  +0.38%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.160"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  -1.24%  "
$ws.Range("E11").Value = "  -1.48%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.510"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +1.56%  "
$ws.Range("E13").Value = "  +0.53%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "39.56"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +5.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.739.22"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.25%  "
$ws.Range("E16").Value = "  +4.63%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "66.072.23"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.56%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.219.04"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +1.48%  "
$ws.Range("E19").Value = "  +0.22%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "510.91"
$ws.Range("D20").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "15.47"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +3.97%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.741"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.94%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "15.55"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +0.23%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.12"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.92%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "85.05"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  -0.09%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.38"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +2.79%  "
$ws.Range("E28").Value = "  +2.88%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +3.88%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.88"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.83"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +7.49%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "28.16"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.20%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.23"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +2.61%  "
$ws.Range("E34").Value = "  +0.05%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "6.59"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  +0.02%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "55.04"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  -1.18%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.0907"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "487.80"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +2.66%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.98"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -4.02%  "
$ws.Range("E41").Value = "  +3.04%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.304"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +6.80%  "
$ws.Range("E43").Value = "  +2.55%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.957.45"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  -3.79%  "
$ws.Range("E45").Value = "  +0.54%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₃0643"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  +5.80%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "28.65"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  -2.68%  "
$ws.Range("E49").Value = "  +1.25%  "
$ws.Range("E50").Value = "  +2.82%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "120.55"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +0.02%  "
